$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell A1 keeps its text but becomes vertically centered ---
$ws.Range("A1").VerticalAlignment = -4108   # xlCenter

# --- Data rows: re-ordered serial numbers + newly added names ---
# Row 2: "1." -> "2." / name "Anik"
$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"
$a2.HorizontalAlignment = -4131   # xlLeft
$a2.VerticalAlignment = -4160     # xlTop
$a2.Value = "2."
$ws.Range("B2").Value = "Anik"

# Copy A2's format (text/left/top) onto the rest of the serial-number column
# in one shot so no extra intermediate styles get minted.
$a2.Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)   # xlPasteFormats

# Row 3: "2." -> "1." / name "Aditi"
$ws.Range("A3").Value = "1."
$ws.Range("B3").Value = "Aditi"

# Row 4: "3." -> "4." / name "Darpan"
$ws.Range("A4").Value = "4."
$ws.Range("B4").Value = "Darpan"

# Row 5: "4." -> "3." / name "Arnab"
$ws.Range("A5").Value = "3."
$ws.Range("B5").Value = "Arnab"

# --- New trailing row 6: empty, centered placeholder cell ---
$a6 = $ws.Range("A6")
$a6.NumberFormat = "@"
$a6.HorizontalAlignment = -4108   # xlCenter
$a6.VerticalAlignment = -4160     # xlTop

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 8

# --- Selection moves from B5 to A5 ---
$ws.Range("A5").Select()
